# Applies the documented change:
#  - the "Filmy, ktore..." bullet's text (already complete once its two
#    runs are read together) is consolidated into a single run, and it
#    loses its embedded _GoBack bookmark
#  - the "Wyszukanie wypozyczen ze sklepow..." and "Wyszukanie liczby
#    wypozyczen..." bullets each gain a leading checkmark (the same
#    "Segoe UI Symbol" glyph + two space runs used by every other
#    checked bullet in the list), and the _GoBack bookmark is re-created
#    right before "Wyszukanie liczby wypozyczen..." (after its new
#    checkmark prefix).

$d = $word.ActiveDocument

function Find-ParagraphContaining($needle) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# --- locate the three bullets involved -------------------------------
$pFilmy  = Find-ParagraphContaining "Filmy"

# the checkmark run + the two single-space runs that precede "Filmy"
# are the template we copy onto the other two bullets
$filmyStart = $pFilmy.Range.Start
$checkTemplate = $d.Range($filmyStart, $filmyStart + 3)
$checkLen = $checkTemplate.End - $checkTemplate.Start

# --- remove the old (mid-sentence) _GoBack bookmark -------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- consolidate the "Filmy...filmow." text into a single run ---------
# (leave the checkmark + two leading spaces' runs/formatting untouched)
$pFilmy = Find-ParagraphContaining "Filmy"
$textStart = $pFilmy.Range.Start + $checkLen
$textEnd = $pFilmy.Range.End - 1
$textRange = $d.Range($textStart, $textEnd)
$filmyText = $textRange.Text
$textRange.Delete()
$d.Range($textStart, $textStart).InsertAfter($filmyText)

# --- prefix "Wyszukanie liczby wypozyczen..." with the checkmark ------
# (done before the city bullet so paragraph/character offsets for the
# still-untouched "Wyszukanie liczby..." bullet don't shift beneath us)
$pYear = Find-ParagraphContaining "od roku"
$yearInsertPoint = $pYear.Range.Start
$yearDest = $d.Range($yearInsertPoint, $yearInsertPoint)
$yearDest.FormattedText = $checkTemplate.FormattedText

# re-add the bookmark right after the freshly inserted checkmark/spaces,
# i.e. right before "Wyszukanie liczby wypozyczen..."
$bookmarkPos = $yearInsertPoint + $checkLen
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# --- prefix "Wyszukanie wypozyczen ze sklepow..." with the checkmark --
$pCity = Find-ParagraphContaining "konkretnego miasta"
$cityInsertPoint = $pCity.Range.Start
$cityDest = $d.Range($cityInsertPoint, $cityInsertPoint)
$cityDest.FormattedText = $checkTemplate.FormattedText

Write-Host "Done"
